# Refresh the crypto price/volume snapshot (columns D and E) to match
# the latest scrape, row by row, as published by the GitHub Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''69.179.24'
$ws.Range("E2").Value = '  +2.42%  '
$ws.Range("D3").Value = '''3.739.87'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''602.59'
$ws.Range("E5").Value = '  +0.76%  '
$ws.Range("D6").Value = '''168.24'
$ws.Range("E6").Value = '  +0.85%  '
$ws.Range("D7").Value = '''3.738.79'
$ws.Range("E7").Value = '  +0.81%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  +0.89%  '
$ws.Range("D10").Value = '''0.167'
$ws.Range("E10").Value = '  +1.31%  '
$ws.Range("D11").Value = '''6.41'
$ws.Range("E11").Value = '  +3.76%  '
$ws.Range("D12").Value = '''0.461'
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("D13").Value = '''38.12'
$ws.Range("E13").Value = '  +0.77%  '
$ws.Range("D14").Value = '''0.0000248'
$ws.Range("E14").Value = '  +1.97%  '
$ws.Range("D15").Value = '''4.361.30'
$ws.Range("E15").Value = '  +0.59%  '
$ws.Range("D16").Value = '''3.764.25'
$ws.Range("E16").Value = '  +1.34%  '
$ws.Range("D17").Value = '''69.174.86'
$ws.Range("E17").Value = '  +2.36%  '
$ws.Range("D18").Value = '''7.30'
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("E19").Value = '  -0.79%  '
$ws.Range("D20").Value = '''17.03'
$ws.Range("E20").Value = '  -3.79%  '
$ws.Range("D21").Value = '''10.91'
$ws.Range("E21").Value = '  +17.60%  '
$ws.Range("D22").Value = '''494.54'
$ws.Range("E22").Value = '  +1.53%  '
$ws.Range("D23").Value = '''0.726'
$ws.Range("E23").Value = '  -0.28%  '
$ws.Range("D24").Value = '''0.0000152'
$ws.Range("E24").Value = '  +9.38%  '
$ws.Range("D25").Value = '''84.77'
$ws.Range("E25").Value = '  -0.40%  '
$ws.Range("E26").Value = '  +0.59%  '
$ws.Range("E27").Value = '  +0.75%  '
$ws.Range("D28").Value = '''10.16'
$ws.Range("E28").Value = '  +0.82%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("E30").Value = '  +1.92%  '
$ws.Range("D31").Value = '''2.51'
$ws.Range("E31").Value = '  +6.36%  '
$ws.Range("D32").Value = '''8.09'
$ws.Range("E32").Value = '  +5.01%  '
$ws.Range("D33").Value = '''31.60'
$ws.Range("E33").Value = '  +0.55%  '
$ws.Range("D34").Value = '''3.882.84'
$ws.Range("E34").Value = '  +0.73%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").Value = '''3.671.82'
$ws.Range("E36").Value = '  +0.41%  '
$ws.Range("D37").Value = '''0.999'
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("D39").Value = '''5.88'
$ws.Range("E39").Value = '  +1.01%  '
$ws.Range("E40").Value = '  +2.34%  '
$ws.Range("D41").Value = '''0.324'
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("D42").Value = '''3.01'
$ws.Range("E42").Value = '  +6.91%  '
$ws.Range("D43").Value = '''432.86'
$ws.Range("E43").Value = '  +1.35%  '
$ws.Range("E44").Value = '  -0.35%  '
$ws.Range("D45").Value = '''1.99'
$ws.Range("E45").Value = '  +2.64%  '
$ws.Range("D46").Value = '''8.49'
$ws.Range("E46").Value = '  +0.32%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").Value = '''40.31'
$ws.Range("E48").Value = '  -0.30%  '
$ws.Range("D49").Value = '''140.58'
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("D50").Value = '''2.778.60'
$ws.Range("E50").Value = '  +1.31%  '
$ws.Range("E51").Value = '  +0.40%  '
